# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-29, replacing the previous
# "Strike#"-derived values with the recalculated K values.
$kValues = @(2,1,3,1,1,1,0,1,1,1,2,4,2,3,0,0,1,2,2,1,1,1,0,0,0,2,2,1)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
